$wb = $excel.ActiveWorkbook

# --- Add the new "testbox" worksheet as the last sheet (after OneJN_TC1) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "testbox"

# --- Row 1: header labels (A1:D1 first, then E1:H1) ---
$newSheet.Range("A1").Value = "fullname"
$newSheet.Range("B1").Value = "email"
$newSheet.Range("C1").Value = "curaddress"
$newSheet.Range("D1").Value = "peraddress"

# --- Row 2: data values, entered in the order A2, C2, B2, D2 ---
$newSheet.Range("A2").Value = "Karun"
$newSheet.Range("C2").Value = "address"
$newSheet.Range("B2").Value = "singh@gmail.com"
$newSheet.Range("D2").Value = "peraddress"

# --- Row 1 continued: expected-value headers ---
$newSheet.Range("E1").Value = "expfullname"
$newSheet.Range("F1").Value = "expemail"
$newSheet.Range("G1").Value = "expcuraddress"
$newSheet.Range("H1").Value = "expperaddress"

# --- Row 2 continued: expected values ---
$newSheet.Range("E2").Value = "Name:Karun"
$newSheet.Range("F2").Value = "Email:singh@gmail.com"
$newSheet.Range("G2").Value = "Current Address :address"
$newSheet.Range("H2").Value = "Permananet Address :peraddress"

# --- Hyperlink the email cell, then restore the workbook's shared Hyperlink style ---
$newSheet.Hyperlinks.Add($newSheet.Range("B2"), "mailto:singh@gmail.com")
$newSheet.Range("B2").Style = "Hyperlink"

# --- Selection on the new sheet ---
$newSheet.Range("I16").Select()
